# Update benchmark: 2026-02-01 06:58:47 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# YKB (column F) fee figures newly populated for EFT rows
$ws.Range("F3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("F5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# YKB (column F) fee figures newly populated for HAVALE rows
$ws.Range("F8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("F10").Value = "15,23 TL - 30,47 TL - 304,71 TL"

# GELEN SWIFT row: İŞBANKASI max updated, YKB value newly populated
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 11.380 TL"
$ws.Range("F13").Value = "Hesaba: Asgari 300 TL | Azami 3.080 TL"

# GİDEN SWIFT - Mobil row: YKB value newly populated
$ws.Range("F14").Value = "1.952,38 TL - 9.523,81 TL"
